$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime M-column cells with the existing date style (s="2") by copying format from M176
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(177, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(178, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(179, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(180, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(181, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(182, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(183, 13))
$ws.Cells.Item(176, 13).Copy($ws.Cells.Item(184, 13))

# Row 177
$ws.Cells.Item(177, 1).Value = 45
$ws.Cells.Item(177, 2).Value = "Male"
$ws.Cells.Item(177, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(177, 4).Value = "Distance"
$ws.Cells.Item(177, 5).Value = "20/30"
$ws.Cells.Item(177, 6).Value = "20/50"
$ws.Cells.Item(177, 7).Value = "20/100"
$ws.Cells.Item(177, 8).Value = "20/70"
$ws.Cells.Item(177, 9).Value = "1-2 years"
$ws.Cells.Item(177, 10).Value = "Uninsured"
$ws.Cells.Item(177, 11).Value = "Cost"
$ws.Cells.Item(177, 12).Value = "Y"
$ws.Cells.Item(177, 13).Value = 45318
$ws.Cells.Item(177, 14).Value = "Dia De La Mujer Latina"

# Row 178
$ws.Cells.Item(178, 1).Value = 36
$ws.Cells.Item(178, 2).Value = "Male"
$ws.Cells.Item(178, 3).Value = "Multiple"
$ws.Cells.Item(178, 4).Value = "Both"
$ws.Cells.Item(178, 5).Value = "20/70"
$ws.Cells.Item(178, 6).Value = "20/50"
$ws.Cells.Item(178, 7).Value = "20/40"
$ws.Cells.Item(178, 8).Value = "20/30"
$ws.Cells.Item(178, 9).Value = "3-5 years"
$ws.Cells.Item(178, 10).Value = "Uninsured"
$ws.Cells.Item(178, 11).Value = "Cost, Lack of Knowledge"
$ws.Cells.Item(178, 12).Value = "Y"
$ws.Cells.Item(178, 13).Value = 45318
$ws.Cells.Item(178, 14).Value = "Dia De La Mujer Latina"

# Row 179
$ws.Cells.Item(179, 1).Value = 41
$ws.Cells.Item(179, 2).Value = "Female"
$ws.Cells.Item(179, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(179, 5).Value = "20/40"
$ws.Cells.Item(179, 6).Value = "20/50"
$ws.Cells.Item(179, 7).Value = "unable to screen"
$ws.Cells.Item(179, 8).Value = "unable to screen"
$ws.Cells.Item(179, 9).Value = "Never"
$ws.Cells.Item(179, 10).Value = "Uninsured"
$ws.Cells.Item(179, 11).Value = "Cost, Clinic Waiting Time, Lack of Knowledge, Other"
$ws.Cells.Item(179, 12).Value = "Y"
$ws.Cells.Item(179, 13).Value = 45318
$ws.Cells.Item(179, 14).Value = "Dia De La Mujer Latina"

# Row 180
$ws.Cells.Item(180, 1).Value = 47
$ws.Cells.Item(180, 2).Value = "Male"
$ws.Cells.Item(180, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(180, 5).Value = "20/30"
$ws.Cells.Item(180, 6).Value = "20/20"
$ws.Cells.Item(180, 7).Value = "20/100"
$ws.Cells.Item(180, 8).Value = "20/50"
$ws.Cells.Item(180, 9).Value = "1-2 years"
$ws.Cells.Item(180, 11).Value = "Cost"
$ws.Cells.Item(180, 12).Value = "Y"
$ws.Cells.Item(180, 13).Value = 45318
$ws.Cells.Item(180, 14).Value = "Dia De La Mujer Latina"

# Row 181
$ws.Cells.Item(181, 1).Value = 42
$ws.Cells.Item(181, 2).Value = "Female"
$ws.Cells.Item(181, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(181, 5).Value = "20/40"
$ws.Cells.Item(181, 6).Value = "20/50"
$ws.Cells.Item(181, 7).Value = "20/40"
$ws.Cells.Item(181, 8).Value = "20/40"
$ws.Cells.Item(181, 9).Value = "More than 5 years"
$ws.Cells.Item(181, 10).Value = "Uninsured"
$ws.Cells.Item(181, 11).Value = "Cost"
$ws.Cells.Item(181, 12).Value = "Y"
$ws.Cells.Item(181, 13).Value = 45318
$ws.Cells.Item(181, 14).Value = "Dia De La Mujer Latina"

# Row 182
$ws.Cells.Item(182, 1).Value = 66
$ws.Cells.Item(182, 2).Value = "Male"
$ws.Cells.Item(182, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(182, 4).Value = "Near"
$ws.Cells.Item(182, 5).Value = "20/50"
$ws.Cells.Item(182, 6).Value = "20/50"
$ws.Cells.Item(182, 7).Value = "20/70"
$ws.Cells.Item(182, 8).Value = "20/50"
$ws.Cells.Item(182, 9).Value = "3-5 years"
$ws.Cells.Item(182, 10).Value = "Uninsured"
$ws.Cells.Item(182, 11).Value = "Cost, Fear of Doctors"
$ws.Cells.Item(182, 12).Value = "Y"
$ws.Cells.Item(182, 13).Value = 45318
$ws.Cells.Item(182, 14).Value = "Dia De La Mujer Latina"

# Row 183
$ws.Cells.Item(183, 1).Value = 28
$ws.Cells.Item(183, 2).Value = "Female"
$ws.Cells.Item(183, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(183, 5).Value = "20/50"
$ws.Cells.Item(183, 6).Value = "20/50"
$ws.Cells.Item(183, 7).Value = "20/70"
$ws.Cells.Item(183, 8).Value = "20/40"
$ws.Cells.Item(183, 9).Value = "1-2 years"
$ws.Cells.Item(183, 10).Value = "Uninsured"
$ws.Cells.Item(183, 11).Value = "Cost"
$ws.Cells.Item(183, 12).Value = "Y"
$ws.Cells.Item(183, 13).Value = 45318
$ws.Cells.Item(183, 14).Value = "Dia De La Mujer Latina"

# Row 184
$ws.Cells.Item(184, 1).Value = 36
$ws.Cells.Item(184, 2).Value = "Male"
$ws.Cells.Item(184, 3).Value = "Hispanic/Latino"
$ws.Cells.Item(184, 5).Value = "20/40"
$ws.Cells.Item(184, 6).Value = "20/40"
$ws.Cells.Item(184, 7).Value = "20/50"
$ws.Cells.Item(184, 8).Value = "20/50"
$ws.Cells.Item(184, 9).Value = "Never"
$ws.Cells.Item(184, 10).Value = "Uninsured"
$ws.Cells.Item(184, 11).Value = "Cost"
$ws.Cells.Item(184, 12).Value = "Y"
$ws.Cells.Item(184, 13).Value = 45318
$ws.Cells.Item(184, 14).Value = "Dia De La Mujer Latina"

# Re-autofit the data columns to reflect the widened content (approximates
# Excel's own autofit pass after the new rows were typed in)
$ws.Columns.Item(3).ColumnWidth = 21.307291666666668
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 17.307291666666668
$ws.Columns.Item(6).ColumnWidth = 18.592447916666668
$ws.Columns.Item(7).ColumnWidth = 16.877604166666668
$ws.Columns.Item(8).ColumnWidth = 18.307291666666668
$ws.Columns.Item(9).ColumnWidth = 21.022135416666668
$ws.Columns.Item(10).ColumnWidth = 13.877604166666666
$ws.Columns.Item(11).ColumnWidth = 36.451822916666664
$ws.Columns.Item(13).ColumnWidth = 25.307291666666668
$ws.Columns.Item(14).ColumnWidth = 25.166666666666668

# Last cell the user had selected after typing the new data
$ws.Range("I9").Select() | Out-Null

